# EncuestasAndalucia2.xlsx — add 3 new poll rows at the top of the data table
# (rows 3-5), pushing the existing history down by three rows, and update the
# active selection / view accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows above the current row 3 (the most recent poll
# entries are stored with the newest at the top, right below the fixed
# header/most-recent row 2). This shifts every existing row from 3..195
# down to 6..198 and picks up the same row styling as row 3 did.
$ws.Range("A3:L5").EntireRow.Insert()

# Row 3: new SigmaDos poll
$ws.Range("A3").Value = 43420
$ws.Range("B3").Value = 43424
$ws.Range("C3").Value = 2500
$ws.Range("D3").Value = "SigmaDos"
$ws.Range("E3").Value = 19
$ws.Range("F3").Value = 30.3
$ws.Range("G3").Value = 20.2
$ws.Range("H3").Value = 20.8
$ws.Range("I3").Value = 5.9

# Row 4: new 40dB poll
$ws.Range("A4").Value = 43416
$ws.Range("B4").Value = 43423
$ws.Range("C4").Value = 1204
$ws.Range("D4").Value = "40dB"
$ws.Range("E4").Value = 19
$ws.Range("F4").Value = 32.1
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = 4.3

# Row 5: new Celeste-Tel poll
$ws.Range("A5").Value = 43410
$ws.Range("B5").Value = 43423
$ws.Range("C5").Value = 2400
$ws.Range("D5").Value = "Celeste-Tel"
$ws.Range("E5").Value = 21.1
$ws.Range("F5").Value = 35.9
$ws.Range("G5").Value = 14.6
$ws.Range("H5").Value = 21.8
$ws.Range("I5").Value = 3.1

# Update the saved view state to match: active cell I5 selected.
$ws.Range("I5").Select()
